# Daily attendance processing - swap the order of "Recorded By" names in
# column G from "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = [string]$cell.Text
    if ($text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
